$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.004.71'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.672.60'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  +0.08%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '214.93'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +1.87%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '21.42'
$cell.Style = "Normal"
$ws.Range('E9').Value = '  +5.37%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '0.0622'
$cell.Style = "Normal"
$ws.Range('E10').Value = '  -0.08%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '0.0888'
$cell.Style = "Normal"
$ws.Range('E11').Value = '  -0.35%  '
$ws.Range('D12').Value = '1.910.20'
$ws.Range('E12').Value = '  +0.35%  '
$ws.Range('D13').Value = '1.713.30'
$ws.Range('E13').Value = '  +2.69%  '
$ws.Range('E14').Value = '  +0.84%  '
$ws.Range('E15').Value = '  +1.53%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '66.22'
$cell.Style = "Normal"
$ws.Range('E16').Value = '  +0.83%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '8.20'
$cell.Style = "Normal"
$ws.Range('E17').Value = '  +2.95%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '27.000.03'
$ws.Range('E18').Value = '  +0.36%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '235.33'
$cell.Style = "Normal"
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').Value = '0.0₃0736'
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').Value = '  +1.74%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '9.26'
$cell.Style = "Normal"
$ws.Range('E23').Value = '  +0.95%  '
$ws.Range('E24').Value = '  -2.20%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '147.00'
$cell.Style = "Normal"
$ws.Range('E25').Value = '  +0.28%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '7.25'
$cell.Style = "Normal"
$ws.Range('E26').Value = '  +1.90%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '16.44'
$cell.Style = "Normal"
$ws.Range('E27').Value = '  +3.57%  '
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('E30').Value = '  +0.71%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('D33').Value = '1.540.11'
$ws.Range('E33').Value = '  +6.41%  '
$ws.Range('E34').Value = '  +0.79%  '
$ws.Range('E35').Value = '  +4.82%  '
$ws.Range('E36').Value = '  -1.11%  '
$ws.Range('E37').Value = '  +0.97%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '0.0174'
$cell.Style = "Normal"
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '0.911'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  +0.88%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '1.04'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  +4.76%  '
$ws.Range('E41').Value = '  +0.00%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '67.64'
$cell.Style = "Normal"
$ws.Range('E42').Value = '  +2.56%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '5.54'
$cell.Style = "Normal"
$ws.Range('E43').Value = '  -3.41%  '
$ws.Range('E44').Value = '  -2.44%  '
$ws.Range('D45').Value = '1.816.71'
$ws.Range('E45').Value = '  +0.70%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '0.781'
$cell.Style = "Normal"
$ws.Range('E46').Value = '  -0.04%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '90.53'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  -0.26%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0106'
$ws.Range('E48').Value = '  +2.41%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '1.54'
$cell.Style = "Normal"
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '0.104'
$cell.Style = "Normal"
$ws.Range('E50').Value = '  +1.89%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '8.03'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  +6.36%  '
